$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 4, pushing the existing data (old rows 4-24)
# down to rows 6-26. Excel copies formatting (incl. the date style on column D)
# from the row above on insert.
$ws.Rows.Item(4).EntireRow.Insert()
$ws.Rows.Item(4).EntireRow.Insert()

# New row 4: Primera, fecha serial 44616
$ws.Range("A4").Value = 11
$ws.Range("B4").Value = "Vega Monumental Concepción"
$ws.Range("C4").Value = "Bíobío"
$ws.Range("D4").Value = 44616
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 100112037
$ws.Range("G4").Value = "Cebollín"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 600
$ws.Range("L4").Value = 700
$ws.Range("M4").Value = 650
$ws.Range("N4").Value = "`$/paquete 6 unidades"
$ws.Range("O4").Value = "Región de Ñuble"
$ws.Range("P4").Value = 108
$ws.Range("Q4").Value = 6
$ws.Range("R4").Value = "Hortaliza"

# New row 5: Segunda, same fecha
$ws.Range("A5").Value = 11
$ws.Range("B5").Value = "Vega Monumental Concepción"
$ws.Range("C5").Value = "Bíobío"
$ws.Range("D5").Value = 44616
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 100112037
$ws.Range("G5").Value = "Cebollín"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 500
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = 500
$ws.Range("N5").Value = "`$/paquete 6 unidades"
$ws.Range("O5").Value = "Región de Ñuble"
$ws.Range("P5").Value = 83
$ws.Range("Q5").Value = 6
$ws.Range("R5").Value = "Hortaliza"
